# Add 2022-Q4 data sheet and update the 总计 (Total) summary sheet.
#
# Net structural effect (per the target diff):
#   - A brand-new "2022-Q4" sheet is inserted right after "总计", in position 2.
#   - Every other quarter sheet ("2022-Q3", "2022-Q1", "2021-Q4", ... "2020-Q4")
#     keeps its own name + data, but shifts one tab position to the right
#     to make room.
#   - The "总计" (Total) sheet gets one new row (2022-Q4 summary) inserted
#     right under its header, and the existing rows shift down with their
#     running index (column A) renumbered to stay sequential.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)         # "总计"
$wsQ3 = $wb.Worksheets.Item(2)        # current "2022-Q3" (becomes the style/template donor)

# ---------------------------------------------------------------------
# 1) Create the new "2022-Q4" sheet by duplicating the "2022-Q3" sheet
#    (right after "总计"). This keeps identical column styles/formats,
#    and bumps every later sheet one slot to the right automatically.
# ---------------------------------------------------------------------
$wsQ3.Copy($null, $ws1)

$wsQ4 = $wb.Worksheets.Item(2)
$wsQ4.Name = "2022-Q4"

# Helper: write a text value into a cell while forcing text storage
# (leading apostrophe), then strip the leftover number-format style so
# the cell ends up with default styling — matching cells that were never
# explicitly formatted.
function Set-TextCell($cell, $text) {
    $cell.Value = "'" + $text
    $cell.ClearFormats()
}

# Row 2 (fund rank 0): code/name unchanged, metrics updated, rank 6 -> 4
Set-TextCell $wsQ4.Cells.Item(2,2) "008115"
Set-TextCell $wsQ4.Cells.Item(2,3) "天弘中证红利低波动100指数C"
Set-TextCell $wsQ4.Cells.Item(2,4) "2.67"
Set-TextCell $wsQ4.Cells.Item(2,5) "94.95"
Set-TextCell $wsQ4.Cells.Item(2,6) "2.02"
Set-TextCell $wsQ4.Cells.Item(2,7) "0.0539"
$wsQ4.Cells.Item(2,8).Value = 4

# Row 3 (fund rank 1): now 515100 / 景顺长城中证红利低波动100ETF
Set-TextCell $wsQ4.Cells.Item(3,2) "515100"
Set-TextCell $wsQ4.Cells.Item(3,3) "景顺长城中证红利低波动100ETF"
Set-TextCell $wsQ4.Cells.Item(3,4) "1.96"
Set-TextCell $wsQ4.Cells.Item(3,5) "99.17"
Set-TextCell $wsQ4.Cells.Item(3,6) "2.12"
Set-TextCell $wsQ4.Cells.Item(3,7) "0.0416"
$wsQ4.Cells.Item(3,8).Value = 4

# Row 4 (fund rank 2): now 008114 / 天弘中证红利低波动100指数A
Set-TextCell $wsQ4.Cells.Item(4,2) "008114"
Set-TextCell $wsQ4.Cells.Item(4,3) "天弘中证红利低波动100指数A"
Set-TextCell $wsQ4.Cells.Item(4,4) "1.98"
Set-TextCell $wsQ4.Cells.Item(4,5) "94.95"
Set-TextCell $wsQ4.Cells.Item(4,6) "2.02"
Set-TextCell $wsQ4.Cells.Item(4,7) "0.0400"
$wsQ4.Cells.Item(4,8).Value = 4

# ---------------------------------------------------------------------
# 2) Update the "总计" (Total) sheet: insert the 2022-Q4 summary row
#    right after the header, then renumber the running index.
# ---------------------------------------------------------------------
$ws1.Rows.Item(2).Insert()

$newRow = $ws1.Range("A2:D2")
$newRow.ClearFormats()

$ws1.Cells.Item(2,1).Value = 0
$ws1.Cells.Item(2,2).Value = "2022-Q4"
$ws1.Cells.Item(2,3).Value = 3
$ws1.Cells.Item(2,4).Value = 0.14

# Restore the running-index column style (bold/border/center) on the new A2,
# copying it from A3 (the old A2, which already carries that style).
$ws1.Cells.Item(3,1).Copy()
$ws1.Cells.Item(2,1).PasteSpecial(-4122)

# Renumber column A (0-based running index) for every row that shifted down.
for ($r = 3; $r -le 9; $r++) {
    $ws1.Cells.Item($r,1).Value = $r - 2
}

# Restore "总计" as the active sheet/selection (matches the original workbook state).
$ws1.Activate()
$ws1.Range("A1").Select()
